{"js": "// 1) Text edit: remove the word \"comercial \" from the \"Oportunidades a serem\n//    exploradas\" paragraph so that\n//    \"...necessidade de investimentos e mesmo comercial de informa\u00e7\u00e3o...\"\n//    becomes\n//    \"...necessidade de investimentos e mesmo de informa\u00e7\u00e3o...\"\nconst target = context.document.body.search(\"mesmo comercial de informa\u00e7\u00e3o\", { matchCase: true });\ntarget.load(\"items,text\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  target.items[0].insertText(\"mesmo de informa\u00e7\u00e3o\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Remove the whole \"Vocabul\u00e1rio de Neg\u00f3cios:\" section (the heading plus its\n//    four bullet paragraphs: Commodities, Ativos, B2B, B2C). The blank\n//    paragraph that used to sit right before the heading is left untouched.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Vocabul\u00e1rio de Neg\u00f3cios\") !== -1) {\n    headingIndex = i;\n    break;\n  }\n}\n\nif (headingIndex !== -1) {\n  // Delete the heading paragraph and the four bullet paragraphs that follow\n  // it (Commodities, Ativos, B2B, B2C) - five paragraphs in total. Delete\n  // from the bottom up so earlier indices stay valid.\n  for (let i = headingIndex + 4; i >= headingIndex; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Text edit: remove the word \"comercial \" from the \"Oportunidades a serem\n#    exploradas\" paragraph so that\n#    \"...necessidade de investimentos e mesmo comercial de informa\u00e7\u00e3o...\"\n#    becomes\n#    \"...necessidade de investimentos e mesmo de informa\u00e7\u00e3o...\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"mesmo comercial de informa\u00e7\u00e3o\"\n$find.Replacement.Text = \"mesmo de informa\u00e7\u00e3o\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Remove the whole \"Vocabul\u00e1rio de Neg\u00f3cios:\" section (the heading plus its\n#    four bullet paragraphs: Commodities, Ativos, B2B, B2C). The blank\n#    paragraph that used to sit right before the heading is left untouched.\n$paras = $d.Paragraphs\n$headingIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    if ($paras.Item($i).Range.Text -like \"*Vocabul\u00e1rio de Neg\u00f3cios*\") {\n        $headingIndex = $i\n        break\n    }\n}\n\nif ($headingIndex -ne -1) {\n    # Delete the heading paragraph and the four bullet paragraphs that follow\n    # it (Commodities, Ativos, B2B, B2C) - five paragraphs in total. Delete\n    # from the bottom up so earlier indices stay valid.\n    for ($i = ($headingIndex + 4); $i -ge $headingIndex; $i--) {\n        $d.Paragraphs.Item($i).Range.Delete()\n    }\n}\n"}
